$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '35.227.86'
$ws.Range("E2").Value = '  -0.12%  '

# Row 3
$ws.Range("D3").Value = '1.908.05'
$ws.Range("E3").Value = '  +0.16%  '

# Row 4
$ws.Range("E4").Value = '  -0.15%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '253.92'
$ws.Range("E5").Value = '  +2.93%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.698'
$ws.Range("E6").Value = '  +0.81%  '

# Row 7
$ws.Range("E7").Value = '  -0.12%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '41.91'
$ws.Range("E8").Value = '  +0.20%  '

# Row 9
$ws.Range("E9").Value = '  +4.16%  '

# Row 10
$ws.Range("E10").Value = '  -2.22%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0763'
$ws.Range("E11").Value = '  +4.66%  '

# Row 12
$ws.Range("E12").Value = '  -1.04%  '

# Row 13
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '13.26'
$ws.Range("E13").Value = '  +7.12%  '

# Row 14
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '2.185.46'
$ws.Range("E14").Value = '  +0.18%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.746'
$ws.Range("E15").Value = '  +5.62%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.03'
$ws.Range("E16").Value = '  +3.57%  '

# Row 17
$ws.Range("D17").Value = '1.915.08'
$ws.Range("E17").Value = '  +0.37%  '

# Row 18
$ws.Range("D18").Value = '35.221.41'
$ws.Range("E18").Value = '  -0.16%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '74.21'
$ws.Range("E19").Value = '  +2.36%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0842'
$ws.Range("E20").Value = '  +2.17%  '

# Row 21
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '243.81'
$ws.Range("E21").Value = '  +0.93%  '

# Row 22
$ws.Range("B22").Value = 'Avalanche'
$ws.Range("C22").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '13.05'
$ws.Range("E22").Value = '  +3.77%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.10'
$ws.Range("E23").Value = '  +5.17%  '

# Row 24
$ws.Range("E24").Value = '  -0.23%  '

# Row 25
$ws.Range("E25").Value = '  +5.95%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.33'
$ws.Range("E26").Value = '  -0.58%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '166.42'
$ws.Range("E27").Value = '  -2.25%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.67'
$ws.Range("E28").Value = '  +0.95%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.61'
$ws.Range("E29").Value = '  +0.87%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.131'
$ws.Range("E30").Value = '  +0.27%  '

# Row 31
$ws.Range("D31").Value = '4.127.85'
$ws.Range("E31").Value = '  -0.56%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.12'
$ws.Range("E32").Value = '  +18.83%  '

# Row 33
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0610'
$ws.Range("E33").Value = '  +6.85%  '

# Row 34
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.33'
$ws.Range("E34").Value = '  +3.22%  '

# Row 35
$ws.Range("E35").Value = '  +19.13%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.22'
$ws.Range("E36").Value = '  +2.13%  '

# Row 37
$ws.Range("E37").Value = '  -0.05%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.868'
$ws.Range("E38").Value = '  -11.16%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.03'
$ws.Range("E39").Value = '  -0.66%  '

# Row 40
$ws.Range("B40").Value = 'Aave'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '100.06'
$ws.Range("E40").Value = '  +10.21%  '

# Row 41
$ws.Range("B41").Value = 'InjectiveProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '17.38'
$ws.Range("E41").Value = '  +6.94%  '

# Row 42
$ws.Range("E42").Value = '  +2.85%  '

# Row 43
$ws.Range("E43").Value = '  +1.57%  '

# Row 44
$ws.Range("E44").Value = '  -3.42%  '

# Row 45
$ws.Range("E45").Value = '  +2.47%  '

# Row 46
$ws.Range("D46").Value = '1.341.14'
$ws.Range("E46").Value = '  -0.40%  '

# Row 47
$ws.Range("E47").Value = '  +0.68%  '

# Row 48
$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.67'
$ws.Range("E48").Value = '  +1.66%  '

# Row 49
$ws.Range("B49").Value = 'MXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.76'
$ws.Range("E49").Value = '  -1.28%  '

# Row 50
$ws.Range("E50").Value = '  -8.33%  '

# Row 51
$ws.Range("B51").Value = 'MultiversX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '43.46'
$ws.Range("E51").Value = '  -8.55%  '
